$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 33
$wsALC.Range("H33").Value = 609.4
$wsALC.Range("I33").Value = 121.333336
$wsALC.Range("K33").Value = 121.333336
$wsALC.Range("M33").Value = 107.666664

# ALC row 103
$wsALC.Range("H103").Value = 795
$wsALC.Range("I103").Value = 790
$wsALC.Range("K103").Value = 2370
$wsALC.Range("M103").Value = -1784

# ARM row 32
$wsARM.Range("H32").Value = 1981.7222
$wsARM.Range("I32").Value = 1981.4667
$wsARM.Range("K32").Value = 1981.4667
$wsARM.Range("M32").Value = -1694.4667

# ARM row 74
$wsARM.Range("H74").Value = 1752.5
$wsARM.Range("I74").Value = 966.3333
$wsARM.Range("J74").Value = 4111
$wsARM.Range("K74").Value = 966.3333
$wsARM.Range("L74").Value = 4111
$wsARM.Range("M74").Value = -92.33330000000001
$wsARM.Range("N74").Value = -5859

# ARM row 77
$wsARM.Range("H77").Value = 1752.5
$wsARM.Range("I77").Value = 966.3333
$wsARM.Range("J77").Value = 4111
$wsARM.Range("K77").Value = 4831.6665
$wsARM.Range("L77").Value = 20555
$wsARM.Range("M77").Value = -463.6665000000003
$wsARM.Range("N77").Value = -29291

# ARM row 110
$wsARM.Range("H110").Value = 1011
$wsARM.Range("I110").Value = 1011
$wsARM.Range("J110").Value = 0
$wsARM.Range("K110").Value = 1011
$wsARM.Range("L110").Value = 0
$wsARM.Range("M110").Value = 1034
$wsARM.Range("N110").ClearContents()

# ARM row 122
$wsARM.Range("H122").Value = 785.8333
$wsARM.Range("I122").Value = 533.6667
$wsARM.Range("J122").Value = 1038
$wsARM.Range("K122").Value = 1601.0001
$wsARM.Range("L122").Value = 3114
$wsARM.Range("M122").Value = 848.9999
$wsARM.Range("N122").Value = -8014

# ARM row 132
$wsARM.Range("H132").Value = 3502.25
$wsARM.Range("I132").Value = 2504.5
$wsARM.Range("J132").Value = 4500
$wsARM.Range("K132").Value = 7513.5
$wsARM.Range("L132").Value = 13500
$wsARM.Range("M132").Value = -4983.5
$wsARM.Range("N132").Value = -18560

# BSM row 24
$wsBSM.Range("H24").Value = 1000
$wsBSM.Range("J24").Value = 1000
$wsBSM.Range("L24").Value = 1000
$wsBSM.Range("N24").Value = -1470

# BSM row 134
$wsBSM.Range("H134").Value = 0
$wsBSM.Range("I134").Value = 0
$wsBSM.Range("K134").Value = 0
$wsBSM.Range("M134").ClearContents()

# CRP row 132
$wsCRP.Range("H132").Value = 5032.8
$wsCRP.Range("I132").Value = 4866.75
$wsCRP.Range("K132").Value = 14600.25
$wsCRP.Range("M132").Value = -12070.25

# CUL row 14
$wsCUL.Range("H14").Value = 665.4286
$wsCUL.Range("I14").Value = 665.4286
$wsCUL.Range("K14").Value = 1996.2858
$wsCUL.Range("M14").Value = -1823.2858

# CUL row 68
$wsCUL.Range("H68").Value = 1853.7142
$wsCUL.Range("I68").Value = 700
$wsCUL.Range("J68").Value = 2046
$wsCUL.Range("K68").Value = 2100
$wsCUL.Range("L68").Value = 6138
$wsCUL.Range("M68").Value = -1289
$wsCUL.Range("N68").Value = -7760

# CUL row 71
$wsCUL.Range("H71").Value = 1853.7142
$wsCUL.Range("I71").Value = 700
$wsCUL.Range("J71").Value = 2046
$wsCUL.Range("K71").Value = 6300
$wsCUL.Range("L71").Value = 18414
$wsCUL.Range("M71").Value = -2244
$wsCUL.Range("N71").Value = -26526

# GSM row 20
$wsGSM.Range("H20").Value = 18800
$wsGSM.Range("I20").Value = 0
$wsGSM.Range("J20").Value = 18800
$wsGSM.Range("K20").Value = 0
$wsGSM.Range("L20").Value = 18800
$wsGSM.Range("M20").ClearContents()
$wsGSM.Range("N20").Value = -19290

# GSM row 24
$wsGSM.Range("H24").Value = 9959.799999999999
$wsGSM.Range("I24").Value = 9949.75
$wsGSM.Range("J24").Value = 10000
$wsGSM.Range("K24").Value = 9949.75
$wsGSM.Range("L24").Value = 10000
$wsGSM.Range("M24").Value = -9776.75
$wsGSM.Range("N24").Value = -10346

# GSM row 31
$wsGSM.Range("H31").Value = 413.125
$wsGSM.Range("I31").Value = 413.125
$wsGSM.Range("K31").Value = 413.125
$wsGSM.Range("M31").Value = -121.125

# GSM row 37
$wsGSM.Range("H37").Value = 413.125
$wsGSM.Range("I37").Value = 413.125
$wsGSM.Range("K37").Value = 413.125
$wsGSM.Range("M37").Value = -136.125

# GSM row 70
$wsGSM.Range("H70").Value = 100000
$wsGSM.Range("I70").Value = 100000
$wsGSM.Range("K70").Value = 100000
$wsGSM.Range("M70").Value = -99730

# GSM row 73
$wsGSM.Range("H73").Value = 100000
$wsGSM.Range("I73").Value = 100000
$wsGSM.Range("K73").Value = 100000
$wsGSM.Range("M73").Value = -99064

# GSM row 122
$wsGSM.Range("H122").Value = 1572.2858
$wsGSM.Range("I122").Value = 1601.2
$wsGSM.Range("J122").Value = 1500
$wsGSM.Range("K122").Value = 4803.6
$wsGSM.Range("L122").Value = 4500
$wsGSM.Range("M122").Value = -2353.6
$wsGSM.Range("N122").Value = -9400

# GSM row 132
$wsGSM.Range("H132").Value = 5500
$wsGSM.Range("I132").Value = 5500
$wsGSM.Range("J132").Value = 0
$wsGSM.Range("K132").Value = 16500
$wsGSM.Range("L132").Value = 0
$wsGSM.Range("M132").Value = -13970
$wsGSM.Range("N132").ClearContents()

# LTW row 25
$wsLTW.Range("H25").Value = 0
$wsLTW.Range("I25").Value = 0
$wsLTW.Range("J25").Value = 0
$wsLTW.Range("K25").Value = 0
$wsLTW.Range("L25").Value = 0
$wsLTW.Range("M25").ClearContents()
$wsLTW.Range("N25").ClearContents()

# LTW row 46
$wsLTW.Range("H46").Value = 3773.2727
$wsLTW.Range("I46").Value = 2000
$wsLTW.Range("J46").Value = 3950.6
$wsLTW.Range("K46").Value = 2000
$wsLTW.Range("L46").Value = 3950.6
$wsLTW.Range("M46").Value = -1812
$wsLTW.Range("N46").Value = -4326.6

# LTW row 68
$wsLTW.Range("H68").Value = 4944.222
$wsLTW.Range("J68").Value = 12250
$wsLTW.Range("L68").Value = 12250
$wsLTW.Range("N68").Value = -13748

# LTW row 71
$wsLTW.Range("H71").Value = 4944.222
$wsLTW.Range("J71").Value = 12250
$wsLTW.Range("L71").Value = 61250
$wsLTW.Range("N71").Value = -68738

# LTW row 132
$wsLTW.Range("H132").Value = 8330.385
$wsLTW.Range("I132").Value = 8208.637000000001
$wsLTW.Range("J132").Value = 9000
$wsLTW.Range("K132").Value = 24625.911
$wsLTW.Range("L132").Value = 27000
$wsLTW.Range("M132").Value = -22095.911
$wsLTW.Range("N132").Value = -32060

# LTW row 136
$wsLTW.Range("H136").Value = 4813
$wsLTW.Range("I136").Value = 4813
$wsLTW.Range("K136").Value = 14439
$wsLTW.Range("M136").Value = -11889

# WVR row 80
$wsWVR.Range("H80").Value = 0
$wsWVR.Range("J80").Value = 0
$wsWVR.Range("L80").Value = 0
$wsWVR.Range("N80").ClearContents()

# WVR row 83
$wsWVR.Range("H83").Value = 0
$wsWVR.Range("J83").Value = 0
$wsWVR.Range("L83").Value = 0
$wsWVR.Range("N83").ClearContents()

# WVR row 132
$wsWVR.Range("H132").Value = 3275
$wsWVR.Range("I132").Value = 3230
$wsWVR.Range("K132").Value = 9690
$wsWVR.Range("M132").Value = -7160

# WVR row 136
$wsWVR.Range("H136").Value = 6499.8335
$wsWVR.Range("I136").Value = 5000
$wsWVR.Range("K136").Value = 15000
$wsWVR.Range("M136").Value = -12450
